# Trade #87 closed at 2026-02-17 21:18:21 - unknown UNKNOWN +0.000%
#
# Updates the live trading results workbook:
#   - Summary sheet roll-up numbers
#   - Strategy Status row for MarketMaking
#   - "All Trades" + "MarketMaking" sheets: close out the open trade
#     (Trade #115 / row 116 in "All Trades", row 83 in "MarketMaking")
#     and append the newly opened trade (Trade #148).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Summary
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1401.07   # Current Capital
$summary.Range("B4").Value = 0.86      # Total P&L $
$summary.Range("B5").Value = 0.15      # Total P&L %
$summary.Range("B6").Value = 115       # Total Trades
$summary.Range("B8").Value = 44        # Losing Trades
$summary.Range("B9").Value = 44.35     # Win Rate %

# ---------------------------------------------------------------------
# Strategy Status - MarketMaking row (row 5)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C5").Value = 101.07     # Capital
$status.Range("D5").Value = 82         # Trades
$status.Range("E5").Value = 0.75       # P&L $
$status.Range("F5").Value = 1.07       # P&L %
$status.Range("G5").Value = 45.12      # Win Rate %

# ---------------------------------------------------------------------
# All Trades - close Trade #115 (row 116)
# ---------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")
$allTrades.Range("G116").Value = 0.93
$allTrades.Range("H116").Value = "CLOSED"
$allTrades.Range("I116").Value = -3.125
$allTrades.Range("J116").Value = -0.03
$allTrades.Range("K116").Value = 101.07
$allTrades.Range("L116").Value = "early_exit"
$allTrades.Range("M116").Value = 0.14

# ---------------------------------------------------------------------
# All Trades - append new Trade #148 (row 149)
# ---------------------------------------------------------------------
$allTrades.Range("A149").Value = 148
$allTrades.Range("B149").NumberFormat = "@"
$allTrades.Range("C149").NumberFormat = "@"
$allTrades.Range("B149").Value = "2026-02-17"
$allTrades.Range("C149").Value = "21:18:14"
$allTrades.Range("D149").Value = "MarketMaking"
$allTrades.Range("E149").Value = "DOWN"
$allTrades.Range("F149").Value = 0.96
$allTrades.Range("H149").Value = "OPEN"
$allTrades.Range("I149").Value = 0
$allTrades.Range("J149").Value = 0
$allTrades.Range("K149").Value = 101.0996151053151
$allTrades.Range("M149").Value = 0
$allTrades.Range("N149").Value = 0
$allTrades.Range("O149").Value = 0
$allTrades.Range("P149").Value = 0.6
$allTrades.Range("Q149").Value = "Normal spread capture: 19600 bps"

# ---------------------------------------------------------------------
# MarketMaking - close Trade #115 (row 83)
# ---------------------------------------------------------------------
$mm = $wb.Worksheets.Item("MarketMaking")
$mm.Range("G83").Value = 0.93
$mm.Range("H83").Value = "CLOSED"
$mm.Range("I83").Value = -3.125
$mm.Range("J83").Value = -0.03
$mm.Range("K83").Value = 101.07
$mm.Range("P83").Value = "early_exit"
$mm.Range("Q83").Value = 0.14

# ---------------------------------------------------------------------
# MarketMaking - append new Trade #148 (row 116)
# ---------------------------------------------------------------------
$mm.Range("A116").Value = 148
$mm.Range("B116").NumberFormat = "@"
$mm.Range("C116").NumberFormat = "@"
$mm.Range("B116").Value = "2026-02-17"
$mm.Range("C116").Value = "21:18:14"
$mm.Range("D116").Value = "MarketMaking"
$mm.Range("E116").Value = "DOWN"
$mm.Range("F116").Value = 0.96
$mm.Range("H116").Value = "OPEN"
$mm.Range("I116").Value = 0
$mm.Range("J116").Value = 0
$mm.Range("K116").Value = 101.0996151053151
$mm.Range("L116").Value = 0
$mm.Range("M116").Value = 0
$mm.Range("N116").Value = 0.6
$mm.Range("O116").Value = "Normal spread capture: 19600 bps"
$mm.Range("Q116").Value = 0
